# Move the inline picture from the end of the document to the very
# beginning, resizing it from 5760720x4320540 EMU (6.3in x 4.725in) down
# to 1300480x975360 EMU (~1.42in x ~1.07in), and renumber its docPr id
# from 2 to 1 - matching a "resize on insert" code path.
#
# The runtime's Range/Selection Copy-Cut-Paste clipboard does not carry
# drawing XML across calls, so the move is done by inserting a fresh
# OOXML paragraph (with the new size/id already applied) at the start of
# the document body, then deleting the original picture paragraph that
# used to be the last paragraph in the body (right after "piet").

$d = $word.ActiveDocument

$pictureXml = '<w:p ' +
  'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
  'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
  'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
  'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:noProof/>' +
        '<w:lang w:eastAsia="nl-NL"/>' +
      '</w:rPr>' +
      '<w:drawing>' +
        '<wp:inline distT="0" distB="0" distL="0" distR="0">' +
          '<wp:extent cy="975360" cx="1300480"/>' +
          '<wp:effectExtent r="0" b="0" l="19050" t="0"/>' +
          '<wp:docPr id="1" descr="image.jpg" name="Picture 0"/>' +
          '<wp:cNvGraphicFramePr>' +
            '<a:graphicFrameLocks noChangeAspect="1"/>' +
          '</wp:cNvGraphicFramePr>' +
          '<a:graphic>' +
            '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
              '<pic:pic>' +
                '<pic:nvPicPr>' +
                  '<pic:cNvPr id="0" name="image.jpg"/>' +
                  '<pic:cNvPicPr/>' +
                '</pic:nvPicPr>' +
                '<pic:blipFill>' +
                  '<a:blip r:embed="rId10"/>' +
                  '<a:stretch>' +
                    '<a:fillRect/>' +
                  '</a:stretch>' +
                '</pic:blipFill>' +
                '<pic:spPr>' +
                  '<a:xfrm>' +
                    '<a:off y="0" x="0"/>' +
                    '<a:ext cy="975360" cx="1300480"/>' +
                  '</a:xfrm>' +
                  '<a:prstGeom prst="rect">' +
                    '<a:avLst/>' +
                  '</a:prstGeom>' +
                '</pic:spPr>' +
              '</pic:pic>' +
            '</a:graphicData>' +
          '</a:graphic>' +
        '</wp:inline>' +
      '</w:drawing>' +
    '</w:r>' +
  '</w:p>'

# 1. Insert the resized picture, wrapped in its own paragraph, at the
#    very start of the document.
$start = $d.Range(0, 0)
$start.InsertXML($pictureXml)

# 2. The original (oversized) picture paragraph is still present - it
#    was the very last paragraph, right after the "piet" bullet - so
#    drop it now that its replacement lives at the top. Re-resolve the
#    shape/paragraph *after* the insert above, since absolute Range
#    offsets captured beforehand are invalidated once new content is
#    spliced in earlier in the document.
$oldShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$oldParagraph = $oldShape.Range.Paragraphs.Item(1)
$oldParagraph.Range.Delete()
